$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 1.129852333333333
$ws.Cells.Item(2, 8).Value = 3.389557
$ws.Cells.Item(2, 9).Value = 0.103499333151224
$ws.Cells.Item(2, 10).Value = 0.103499333151224
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 7.759559666666667
$ws.Cells.Item(2, 14).Value = 23.278679
$ws.Cells.Item(2, 15).Value = 0.1480404846036854
$ws.Cells.Item(2, 16).Value = 0.1480404846036854
$ws.Cells.Item(2, 17).Value = 8.767156595022556
$ws.Cells.Item(2, 18).Value = 78.904409355203
$ws.Cells.Item(2, 19).Value = 0.01532209143586548
$ws.Cells.Item(2, 20).Value = 0.01532209143586548

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 1.129852333333333
$ws.Cells.Item(3, 8).Value = 3.389557
$ws.Cells.Item(3, 9).Value = 0.103499333151224
$ws.Cells.Item(3, 10).Value = 0.103499333151224
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 34.70130533333333
$ws.Cells.Item(3, 14).Value = 104.103916
$ws.Cells.Item(3, 15).Value = 0.6620476262326294
$ws.Cells.Item(3, 16).Value = 0.6620476262326295
$ws.Cells.Item(3, 17).Value = 39.20735080057911
$ws.Cells.Item(3, 18).Value = 352.866157205212
$ws.Cells.Item(3, 19).Value = 0.06852148782942792
$ws.Cells.Item(3, 20).Value = 0.06852148782942792

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 1.129852333333333
$ws.Cells.Item(4, 8).Value = 3.389557
$ws.Cells.Item(4, 9).Value = 0.103499333151224
$ws.Cells.Item(4, 10).Value = 0.103499333151224
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 0.375896
$ws.Cells.Item(4, 14).Value = 1.127688
$ws.Cells.Item(4, 15).Value = 0.007171518538563155
$ws.Cells.Item(4, 16).Value = 0.007171518538563155
$ws.Cells.Item(4, 17).Value = 0.4247069726906667
$ws.Cells.Item(4, 18).Value = 3.822362754216
$ws.Cells.Item(4, 19).Value = 0.0007422473864229269
$ws.Cells.Item(4, 20).Value = 0.0007422473864229269

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 1.129852333333333
$ws.Cells.Item(5, 8).Value = 3.389557
$ws.Cells.Item(5, 9).Value = 0.103499333151224
$ws.Cells.Item(5, 10).Value = 0.103499333151224
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 9.578358333333332
$ws.Cells.Item(5, 14).Value = 28.735075
$ws.Cells.Item(5, 15).Value = 0.1827403706251221
$ws.Cells.Item(5, 16).Value = 0.1827403706251221
$ws.Cells.Item(5, 17).Value = 10.82213051241944
$ws.Cells.Item(5, 18).Value = 97.399174611775
$ws.Cells.Item(5, 19).Value = 0.01891350649950765
$ws.Cells.Item(5, 20).Value = 0.01891350649950765

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 2.720340333333333
$ws.Cells.Item(6, 8).Value = 8.161021
$ws.Cells.Item(6, 9).Value = 0.2491948745317264
$ws.Cells.Item(6, 10).Value = 0.2491948745317264
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 7.759559666666667
$ws.Cells.Item(6, 14).Value = 23.278679
$ws.Cells.Item(6, 15).Value = 0.1480404846036854
$ws.Cells.Item(6, 16).Value = 0.1480404846036854
$ws.Cells.Item(6, 17).Value = 21.10864313013989
$ws.Cells.Item(6, 18).Value = 189.977788171259
$ws.Cells.Item(6, 19).Value = 0.03689092998643137
$ws.Cells.Item(6, 20).Value = 0.03689092998643137

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 2.720340333333333
$ws.Cells.Item(7, 8).Value = 8.161021
$ws.Cells.Item(7, 9).Value = 0.2491948745317264
$ws.Cells.Item(7, 10).Value = 0.2491948745317264
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 34.70130533333333
$ws.Cells.Item(7, 14).Value = 104.103916
$ws.Cells.Item(7, 15).Value = 0.6620476262326294
$ws.Cells.Item(7, 16).Value = 0.6620476262326295
$ws.Cells.Item(7, 17).Value = 94.39936051758177
$ws.Cells.Item(7, 18).Value = 849.594244658236
$ws.Cells.Item(7, 19).Value = 0.1649788751530674
$ws.Cells.Item(7, 20).Value = 0.1649788751530674

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 2.720340333333333
$ws.Cells.Item(8, 8).Value = 8.161021
$ws.Cells.Item(8, 9).Value = 0.2491948745317264
$ws.Cells.Item(8, 10).Value = 0.2491948745317264
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 0.375896
$ws.Cells.Item(8, 14).Value = 1.127688
$ws.Cells.Item(8, 15).Value = 0.007171518538563155
$ws.Cells.Item(8, 16).Value = 0.007171518538563155
$ws.Cells.Item(8, 17).Value = 1.022565049938667
$ws.Cells.Item(8, 18).Value = 9.203085449448
$ws.Cells.Item(8, 19).Value = 0.001787105662419195
$ws.Cells.Item(8, 20).Value = 0.001787105662419195

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 2.720340333333333
$ws.Cells.Item(9, 8).Value = 8.161021
$ws.Cells.Item(9, 9).Value = 0.2491948745317264
$ws.Cells.Item(9, 10).Value = 0.2491948745317264
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 9.578358333333332
$ws.Cells.Item(9, 14).Value = 28.735075
$ws.Cells.Item(9, 15).Value = 0.1827403706251221
$ws.Cells.Item(9, 16).Value = 0.1827403706251221
$ws.Cells.Item(9, 17).Value = 26.05639450128611
$ws.Cells.Item(9, 18).Value = 234.507550511575
$ws.Cells.Item(9, 19).Value = 0.04553796372980847
$ws.Cells.Item(9, 20).Value = 0.04553796372980847

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 4.133464333333333
$ws.Cells.Item(10, 8).Value = 12.400393
$ws.Cells.Item(10, 9).Value = 0.3786431106817516
$ws.Cells.Item(10, 10).Value = 0.3786431106817515
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 7.759559666666667
$ws.Cells.Item(10, 14).Value = 23.278679
$ws.Cells.Item(10, 15).Value = 0.1480404846036854
$ws.Cells.Item(10, 16).Value = 0.1480404846036854
$ws.Cells.Item(10, 17).Value = 32.07386312453855
$ws.Cells.Item(10, 18).Value = 288.664768120847
$ws.Cells.Item(10, 19).Value = 0.05605450959717341
$ws.Cells.Item(10, 20).Value = 0.0560545095971734

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 4.133464333333333
$ws.Cells.Item(11, 8).Value = 12.400393
$ws.Cells.Item(11, 9).Value = 0.3786431106817516
$ws.Cells.Item(11, 10).Value = 0.3786431106817515
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 34.70130533333333
$ws.Cells.Item(11, 14).Value = 104.103916
$ws.Cells.Item(11, 15).Value = 0.6620476262326294
$ws.Cells.Item(11, 16).Value = 0.6620476262326295
$ws.Cells.Item(11, 17).Value = 143.4366079154431
$ws.Cells.Item(11, 18).Value = 1290.929471238988
$ws.Cells.Item(11, 19).Value = 0.2506797726161924
$ws.Cells.Item(11, 20).Value = 0.2506797726161924

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 4.133464333333333
$ws.Cells.Item(12, 8).Value = 12.400393
$ws.Cells.Item(12, 9).Value = 0.3786431106817516
$ws.Cells.Item(12, 10).Value = 0.3786431106817515
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 0.375896
$ws.Cells.Item(12, 14).Value = 1.127688
$ws.Cells.Item(12, 15).Value = 0.007171518538563155
$ws.Cells.Item(12, 16).Value = 0.007171518538563155
$ws.Cells.Item(12, 17).Value = 1.553752709042667
$ws.Cells.Item(12, 18).Value = 13.983774381384
$ws.Cells.Item(12, 19).Value = 0.002715446087753402
$ws.Cells.Item(12, 20).Value = 0.002715446087753402

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 4.133464333333333
$ws.Cells.Item(13, 8).Value = 12.400393
$ws.Cells.Item(13, 9).Value = 0.3786431106817516
$ws.Cells.Item(13, 10).Value = 0.3786431106817515
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 9.578358333333332
$ws.Cells.Item(13, 14).Value = 28.735075
$ws.Cells.Item(13, 15).Value = 0.1827403706251221
$ws.Cells.Item(13, 16).Value = 0.1827403706251221
$ws.Cells.Item(13, 17).Value = 39.59180254271944
$ws.Cells.Item(13, 18).Value = 356.326222884475
$ws.Cells.Item(13, 19).Value = 0.0691933823806324
$ws.Cells.Item(13, 20).Value = 0.06919338238063238

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 2.932861
$ws.Cells.Item(14, 8).Value = 8.798583000000001
$ws.Cells.Item(14, 9).Value = 0.2686626816352981
$ws.Cells.Item(14, 10).Value = 0.268662681635298
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 7.759559666666667
$ws.Cells.Item(14, 14).Value = 23.278679
$ws.Cells.Item(14, 15).Value = 0.1480404846036854
$ws.Cells.Item(14, 16).Value = 0.1480404846036854
$ws.Cells.Item(14, 17).Value = 22.75770992353967
$ws.Cells.Item(14, 18).Value = 204.819389311857
$ws.Cells.Item(14, 19).Value = 0.03977295358421518
$ws.Cells.Item(14, 20).Value = 0.03977295358421518

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 2.932861
$ws.Cells.Item(15, 8).Value = 8.798583000000001
$ws.Cells.Item(15, 9).Value = 0.2686626816352981
$ws.Cells.Item(15, 10).Value = 0.268662681635298
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 34.70130533333333
$ws.Cells.Item(15, 14).Value = 104.103916
$ws.Cells.Item(15, 15).Value = 0.6620476262326294
$ws.Cells.Item(15, 16).Value = 0.6620476262326295
$ws.Cells.Item(15, 17).Value = 101.7741050612253
$ws.Cells.Item(15, 18).Value = 915.9669455510281
$ws.Cells.Item(15, 19).Value = 0.1778674906339417
$ws.Cells.Item(15, 20).Value = 0.1778674906339417

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 2.932861
$ws.Cells.Item(16, 8).Value = 8.798583000000001
$ws.Cells.Item(16, 9).Value = 0.2686626816352981
$ws.Cells.Item(16, 10).Value = 0.268662681635298
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 0.375896
$ws.Cells.Item(16, 14).Value = 1.127688
$ws.Cells.Item(16, 15).Value = 0.007171518538563155
$ws.Cells.Item(16, 16).Value = 0.007171518538563155
$ws.Cells.Item(16, 17).Value = 1.102450718456
$ws.Cells.Item(16, 18).Value = 9.922056466104001
$ws.Cells.Item(16, 19).Value = 0.001926719401967631
$ws.Cells.Item(16, 20).Value = 0.001926719401967631

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 2.932861
$ws.Cells.Item(17, 8).Value = 8.798583000000001
$ws.Cells.Item(17, 9).Value = 0.2686626816352981
$ws.Cells.Item(17, 10).Value = 0.268662681635298
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 9.578358333333332
$ws.Cells.Item(17, 14).Value = 28.735075
$ws.Cells.Item(17, 15).Value = 0.1827403706251221
$ws.Cells.Item(17, 16).Value = 0.1827403706251221
$ws.Cells.Item(17, 17).Value = 28.09199359985833
$ws.Cells.Item(17, 18).Value = 252.827942398725
$ws.Cells.Item(17, 19).Value = 0.04909551801517354
$ws.Cells.Item(17, 20).Value = 0.04909551801517353
